$d = $word.ActiveDocument

# ativacao
$d.Content.Find.Execute("Ativação: 01/01/2018", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2021", 2) | Out-Null

# objetivos_pt
$d.Content.Find.Execute("a) Fornecer ao aluno uma visão integrada do processo de desenvolvimento de produto, desde as etapas iniciais de geração da ideia, avaliação econômica e desenvolvimento do conceito do produto até a preparação da fábrica, produção e lançamento do produto.b) Apresentar como os principais conhecimentos da formação de engenheiro de produção podem ser aplicados no processo de desenvolvimento de produto.c) Propiciar uma experiência prática de projeto.", $true, $false, $false, $false, $false, $true, 1, $false, "a) Fornecer ao aluno uma visão integrada do processo de desenvolvimento de produto, desde as etapas iniciais de geração da ideia, avaliação econômica e desenvolvimento do conceito do produto até a preparação da fábrica, produção e lançamento do produto. b) Apresentar como os principais conhecimentos da formação de engenheiro de produção podem ser aplicados no processo de desenvolvimento de produto. c) Propiciar uma experiência prática de projeto.", 2) | Out-Null

# docente
$d.Content.Find.Execute("5840560 - Marco Antonio Carvalho Pereira", $true, $false, $false, $false, $false, $true, 1, $false, "5701460 - Antonio Iacono", 2) | Out-Null

# resumido_pt
$d.Content.Find.Execute("Conceito de processo de negócio e modelo de referência; Desenvolvimento de produto como um processo: Fases principais do processo de desenvolvimento de produto; Definição e conceitos básicos de gerenciamento de projetos; Planejamento estratégico do produto; Planejamento do projeto: projeto informacional, projeto conceitual, projeto detalhado, lançamento e pós-desenvolvimento.", $true, $false, $false, $false, $false, $true, 1, $false, "1. Conceito de processo de negócio e modelo de referência; 2. Desenvolvimento de produto como um processo: Fases principais do processo de desenvolvimento de produto; 3. Definição e conceitos básicos de gerenciamento de projetos; 4. Planejamento estratégico do produto; 5. Planejamento do Projeto; 6. Projeto Informacional; 7. Projeto Conceitual; 8. Projeto Detalhado; 9. Lançamento e Pós-Desenvolvimento. 10. Empreendedorismo e processo de desenvolvimento de produto.", 2) | Out-Null

# resumido_en
$d.Content.Find.Execute("Concept of business and process reference model; Product development as a process: main phases of the product development process; Definition and basic concepts of project management; Strategic product planning; Project planning: information design, conceptual design, detailed design, launch and post-development.", $true, $false, $false, $false, $false, $true, 1, $false, "1.Concept of business and process reference model; 2. Product development as a process: main phases of the product development process; 3. Definition and basic concepts of project management; 4. Strategic product planning; 5. Project planning; 6. information design; 7. conceptual design; 8. detailed design; 9. launch and post-development. 10. Entrepreneurship and product development process", 2) | Out-Null

# programa_pt
$d.Content.Find.Execute("1.Visão geral do processo de desenvolvimento de produto. Conceito de processo de negócio. Modelo de referência e desenvolvimento de produto como um processo. Fases principais do processo de desenvolvimento de produto.Definição e conceitos básicos de gerenciamento de projetos.2.Apresentação de um modelo de referência para desenvolvimento de produto. Descrição do modelo de referência. Importância e uso do modelo de referência. Principais características do modelo de referência para desenvolvimento de produto. 3.Apresentação das fases do processo de desenvolvimento de produto. Nesta parte são apresentados os conceitos e metodologias principais de cada fase do desenvolvimento de produto. Paralelamente os alunos, em grupo, realizam as atividades de desenvolvimento de um produto exemplo conforme tema sugerido pelo professor.3.1.Planejamento estratégico de produtos. Conteúdo sobre a importância e como preparar a estratégia de produto, realizar o gerenciamento da carteira de projetos e elaborar a Minuta do Projeto. 3.2.Planejamento do projeto. Conteúdo sobre técnicas de planejamento do projeto e análise de viabilidade econômica. 3.3.Projeto informacional. Conceitos gerais de pesquisa de mercado, tipos de fontes de dados para pesquisa de mercado, diferenças entre as fontes de dados, formas de utilização de pesquisas e aprofundamento no tipo de pesquisa baseado em enquête. Apresenta-se também em detalhes a técnica do desdobramento da função qualidade (QFD).3.4.Projeto conceitual. Apresentam-se os diversos métodos criativos para concepção, tais como Diagrama de funções, Matriz Morfológica e Matriz de Pugh, e demonstra-se com exemplos como utilizá-los de forma a gerar e escolher uma concepção de produto robusta. Em seguida, os conceitos de Design for X, em especial o Projeto para Montagem, Projeto para Manufatura e Projeto para o Meio-Ambiente (eco-design).3.5.Projeto detalhado e preparação da produção. Apresenta-se o objetivo e resumo das atividades realizadas nas duas fases. Aprofunda-se em aspectos relacionados com a atuação profissional do engenheiro de produção, são eles: classificação de itens, parâmetros críticos, decisão fazer ou comprar e FMEA. No pós-desenvolvimento discute-se especialmente os aspectos relacionados com o meio ambiente, como descarte e reuso de peças do produto.4.Aulas Práticas. Oferecer experiência prática em métodos e etapas do PDP, com o intuito de complementar o aprendizado teórico à partir do desenvolvimento de habilidades e atitudes. As práticas são intercaladas aos módulos de exposição teóricos.", $true, $false, $false, $false, $false, $true, 1, $false, "1. Visão geral do processo de desenvolvimento de produto. Conceito de processo de negócio. Modelo de referência e desenvolvimento de produto como um processo. Fases principais do processo de desenvolvimento de produto. Definição e conceitos básicos de gerenciamento de projetos. 2. Apresentação de um modelo de referência para desenvolvimento de produto. Descrição do modelo de referência. Importância e uso do modelo de referência. Principais características do modelo de referência para desenvolvimento de produto. 3. Apresentação das fases do processo de desenvolvimento de produto. Nesta parte são apresentados os conceitos e metodologias principais de cada fase do desenvolvimento de produto. 3.1. Planejamento estratégico de produtos. Conteúdo sobre a importância e como preparar a estratégia de produto, realizar o gerenciamento da carteira de projetos e elaborar a Minuta do Projeto. 3.2. Planejamento do projeto. Conteúdo sobre técnicas de planejamento do projeto e análise de viabilidade econômica. 3.3. Projeto informacional. Conceitos gerais de pesquisa de mercado, tipos de fontes de dados para pesquisa de mercado, diferenças entre as fontes de dados, formas de utilização de pesquisas e aprofundamento no tipo de pesquisa baseado em enquete. 3.4. Projeto conceitual. Apresentam-se os diversos métodos criativos para concepção, tais como Diagrama de funções, Matriz Morfológica e Matriz de Pugh, e demonstra-se com exemplos como utilizá-los de forma a gerar e escolher uma concepção de produto robusta.3.5. Projeto detalhado e preparação da produção. Apresenta-se o objetivo e resumo das atividades realizadas nas duas fases. No pós-desenvolvimento discute-se especialmente os aspectos relacionados com o meio ambiente, como descarte e reuso de peças do produto. 4. Empreendedorismo e processo de desenvolvimento de produto. O ambiente de inovação e o processo de desenvolvimento de produtos: empresas nascentes de base tecnológica, startups e aceleradoras.", 2) | Out-Null

# programa_en
$d.Content.Find.Execute("1. Overview of the process of product development. Concept of business process. Reference model and product development as a process. Main phases of the product development process. Definition and basic concepts of project management. 2. presentation of a reference model for product development. Description of the reference model. Importance and use of the reference model. Main features of the reference model for product development. 3. presentation of the stages of the product development process. In this section are presented the main concepts and methodologies of each phase of the product development. In parallel with the students, as a group, the development of a product such as theme suggested by the teacher. 3.1. Strategic planning of products. Information about the importance and how to prepare the product strategy, carry out the project portfolio management and draw up the draft of the project. 3.2. Project planning. Information about project planning techniques and analysis of economic feasibility. 3.3. informational Project. General concepts of market research, types of data sources for market research, differences between the data sources, forms of use of research and study on the type of survey-based research. Also presented in detail the technique of quality function deployment (QFD). 3.4. Conceptual design. The various creative methods to design, such as functions, Morphological Matrix diagram and array of Pugh, and demonstrates with examples how to use them to generate and choose a robust product design. Then, the Design concepts for X, in particular the mounting Design, design for manufacturing and design for the environment (eco-design). 3.5. Detailed design and preparation of production. The objective and summary of activities carried out in two phases. Deepened in aspects related to the professional performance of the production engineer, they are: classification of items, critical parameters, make or buy decision and FMEA. In the post-development discusses in particular the aspects relating to the environment, such as disposal and reuse of parts of the product. 4. Practical lessons. Offer hands-on experience in methods and steps of the PDP, in order to complement the theoretical learning from the development of skills and attitudes. The practices are interspersed to theoretical Exposition modules.", $true, $false, $false, $false, $false, $true, 1, $false, "1. overview of the process of product development. concept of business process. reference model and product development as a process. main phases of the product development process. definition and basic concepts of project management. 2. presentation of a reference model for product development. description of the reference model. importance and use of the reference model. main features of the reference model for product development. 3. presentation of the stages of the product development process. in this section are presented the main concepts and methodologies of each phase of the product development. 3.1. strategic planning of products. information about the importance and how to prepare the product strategy, carry out the project portfolio management and draw up the draft of the project. 3.2. project planning. information about project planning techniques and analysis of economic feasibility. 3.3. informational project. general concepts of market research, types of data sources for market research, differences between the data sources, forms of use of research and study on the type of survey-based research. 3.4. conceptual design. the various creative methods to design, such as functions, morphological matrix diagram and array of pugh, and demonstrates with examples how to use them to generate and choose a robust product design. 3.5. detailed design and preparation of production. the objective and summary of activities carried out in two phases. in the post-development discusses in particular the aspects relating to the environment, such as disposal and reuse of parts of the product. 4. entrepreneurship and product development process: innovation environment and product development process: technology-based enterprises, startups and business accelerators", 2) | Out-Null

# metodo
$d.Content.Find.Execute("Aulas expositivas e práticas.", $true, $false, $false, $false, $false, $true, 1, $false, "Provas e trabalho", 2) | Out-Null

# criterio
$d.Content.Find.Execute("Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2", $true, $false, $false, $false, $false, $true, 1, $false, "M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre desenvolvimento de produtosM = Média de aproveitamento do aluno", 2) | Out-Null

# recuperacao
$d.Content.Find.Execute("A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.", 2) | Out-Null
